# Add "Errors" and "Warnings" sheets after "Classes", populate "Warnings"
# with the three skipped-row warning messages, and move the active
# tab/selection state to match the target workbook.

$wb = $excel.ActiveWorkbook
$classesSheet = $wb.Worksheets.Item("Classes")

# New sheets are inserted right after "Classes", in order.
$errorsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $classesSheet)
$errorsSheet.Name = "Errors"

$warningsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $errorsSheet)
$warningsSheet.Name = "Warnings"

# Populate the Warnings sheet with the three "no data" messages (leading
# apostrophe doubled so the literal apostrophe survives Excel's
# force-text-prefix handling of a single leading apostrophe). Restoring the
# "Normal" style afterwards drops the auto-applied quotePrefix cell style so
# the cells keep the default (unstyled) format, same as the source file.
$warningsSheet.Range("A1").Value = "''Sheet ""Classes"" Row: 4 No data found between cells ""A"" and ""D"" Skipping this row',"
$warningsSheet.Range("A2").Value = "''Sheet ""Classes"" Row: 7 No data found between cells ""A"" and ""D"" Skipping this row',"
$warningsSheet.Range("A3").Value = "''Sheet ""Classes"" Row: 9 No data found between cells ""A"" and ""D"" Skipping this row',"
$warningsSheet.Range("A1:A3").Style = "Normal"

# Match the recorded selection on the new active sheet.
$warningsSheet.Range("B7:B8").Select() | Out-Null

# Activate the Warnings tab (last sheet) so it becomes the active/selected tab.
$warningsSheet.Activate()

# Best-effort: move the application window to match the recorded window
# geometry (may be a no-op in this headless host).
$excel.ActiveWindow.Left = 29140
$excel.ActiveWindow.Top = -18660
Write-Output "done"
